$d = $word.ActiveDocument

function Add-Run($para, [string]$text, [bool]$italic) {
    $ip = $para.Range
    $ip.Collapse(0)
    $ip.InsertAfter($text)
    if ($italic) {
        $ip.Font.Italic = $true
    }
}

# Anchor: the current last paragraph is 'BILAGA 1 - Fridlysta arter'.
$anchor = $d.Paragraphs.Last

# --- paragraph 0: style='Heading1' ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -2
$t = 'Kn' + [char]0x00E4 + 'rot ' + [char]0x2013 + ' ekologi samt krav p' + [char]0x00E5 + ' livsmilj' + [char]0x00F6 + 'n'
Add-Run $p $t $false
$anchor = $p

# --- paragraph 1: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = 'Kn' + [char]0x00E4 + 'rot ' + [char]0x00E4 + 'r fridlyst enligt 8 och 15 ' + [char]0x00A7 + [char]0x00A7 + ' artskyddsf' + [char]0x00F6 + 'rordningen och klassad som s' + [char]0x00E5 + 'rbar (VU) enligt r' + [char]0x00F6 + 'dlistan 2020. Kn' + [char]0x00E4 + 'rot ' + [char]0x00E4 + 'r beroende av h' + [char]0x00F6 + 'g och j' + [char]0x00E4 + 'mn luftfuktighet i gamla, ost' + [char]0x00F6 + 'rda skogsmilj' + [char]0x00F6 + 'er och ' + [char]0x00E4 + 'r k' + [char]0x00E4 + 'nslig f' + [char]0x00F6 + 'r snabba f' + [char]0x00F6 + 'r' + [char]0x00E4 + 'ndringar av ljus-/vindf' + [char]0x00F6 + 'rh' + [char]0x00E5 + 'llanden eller uttorkning. P' + [char]0x00E5 + ' grund av ett alltf' + [char]0x00F6 + 'r intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 ' + [char]0x00E5 + 'ren och i framtiden bed' + [char]0x00F6 + 'ms minskningstakten uppg' + [char]0x00E5 + ' till 30 (20-40) %. Till f' + [char]0x00F6 + 'ljd av att arten har en dokumenterat h' + [char]0x00F6 + 'gre minskningstakt if' + [char]0x00F6 + 'rh' + [char]0x00E5 + 'llande till sin generationstid ' + [char]0x00E4 + 'n vad som tidigare varit k' + [char]0x00E4 + 'nt (data fr' + [char]0x00E5 + 'n Riksskogstaxeringen) h' + [char]0x00F6 + 'jdes den till hotkategori s' + [char]0x00E5 + 'rbar (VU) i r' + [char]0x00F6 + 'dlistan 2020 (Artdatabanken, 2021).'
Add-Run $p $t $false
$anchor = $p

# --- paragraph 2: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = 'Samuel Johnsons doktorsavhandling '
Add-Run $p $t $false
$t = [char]0x201C + 'Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation' + [char]0x201C
Add-Run $p $t $true
$t = ' (SLU, Uppsala 2014) visar att det kr' + [char]0x00E4 + 'vs v' + [char]0x00E4 + 'l tilltagna skyddszoner f' + [char]0x00F6 + 'r att kn' + [char]0x00E4 + 'rotens v' + [char]0x00E4 + 'xtplatser inte ska ta skada av skogsbruks' + [char]0x00E5 + 'tg' + [char]0x00E4 + 'rder i intilliggande omr' + [char]0x00E5 + 'den: '
Add-Run $p $t $false
$t = [char]0x201C + 'Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.' + [char]0x201D + ' '
Add-Run $p $t $true
$t = 'Vidare '
Add-Run $p $t $false
$t = [char]0x201C + 'More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).' + [char]0x201D
Add-Run $p $t $true
$anchor = $p

# --- paragraph 3: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = 'Johnsons (2014) rekommendation p' + [char]0x00E5 + ' minst 50 meters breda skyddszoner runt kn' + [char]0x00E4 + 'rotens v' + [char]0x00E4 + 'xtplatser motsvarar en areal p' + [char]0x00E5 + ' 0,78 hektar, vilket ligger i linje med andra studier som gjorts p' + [char]0x00E5 + ' k' + [char]0x00E4 + 'nsliga skogsarter: '
Add-Run $p $t $false
$t = [char]0x201C + 'In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).' + [char]0x201D
Add-Run $p $t $true
$anchor = $p

# --- paragraph 4: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkid' + [char]0x00E9 + 'n kn' + [char]0x00E4 + 'rots skyddsbehov. I uppsatsen ber' + [char]0x00F6 + 'rs problemet med uttorkning f' + [char]0x00F6 + 'r v' + [char]0x00E4 + 'xter, bl.a. f' + [char]0x00F6 + 'r kn' + [char]0x00E4 + 'rot, ett problem som blivit accentuerat p' + [char]0x00E5 + ' grund av den p' + [char]0x00E5 + 'g' + [char]0x00E5 + 'ende klimatf' + [char]0x00F6 + 'r' + [char]0x00E4 + 'ndringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen unders' + [char]0x00F6 + 'ks omr' + [char]0x00E5 + 'den med tre olika avst' + [char]0x00E5 + 'nd fr' + [char]0x00E5 + 'n kalhyggeskant med avseende p' + [char]0x00E5 + ' skydd bl.a. f' + [char]0x00F6 + 'r kn' + [char]0x00E4 + 'rot. Det f' + [char]0x00F6 + 'rsta omr' + [char]0x00E5 + 'det har avst' + [char]0x00E5 + 'nd upp till 20 m fr' + [char]0x00E5 + 'n hyggeskant (Strong edge effect), det andra 20 ' + [char]0x2013 + ' 40 m fr' + [char]0x00E5 + 'n hyggeskant (Weak edge effect) och det tredje avser st' + [char]0x00F6 + 'rre avst' + [char]0x00E5 + 'nd fr' + [char]0x00E5 + 'n hyggeskant, d' + [char]0x00E4 + 'r kanteffekten anses vara f' + [char]0x00F6 + 'rsumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt p' + [char]0x00E5 + ' k' + [char]0x00E4 + 'nsliga och r' + [char]0x00F6 + 'dlistade skogsarter vid de kortare avst' + [char]0x00E5 + 'nden till hyggeskant, medan effekt av uttorkning inte konstaterades p' + [char]0x00E5 + ' st' + [char]0x00F6 + 'rre avst' + [char]0x00E5 + 'nd (Interior). F' + [char]0x00F6 + 'r orkid' + [char]0x00E9 + 'n kn' + [char]0x00E4 + 'rot fann man en rik f' + [char]0x00F6 + 'rekomst (upp till 0,06 dm2/m2) p' + [char]0x00E5 + ' stort avst' + [char]0x00E5 + 'nd fr' + [char]0x00E5 + 'n hyggeskant (Interior), medan f' + [char]0x00F6 + 'rekomsten var liten eller n' + [char]0x00E4 + 'rmast f' + [char]0x00F6 + 'rsumbar i de omr' + [char]0x00E5 + 'den som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet p' + [char]0x00E5 + 'pekar att de allt oftare f' + [char]0x00F6 + 'rekommande torra somrarna ger ytterligare sk' + [char]0x00E4 + 'l att ut' + [char]0x00F6 + 'ka skyddsavst' + [char]0x00E5 + 'ndet fr' + [char]0x00E5 + 'n hyggen till den fuktkr' + [char]0x00E4 + 'vande arten kn' + [char]0x00E4 + 'rot (Koelmeijer m.fl., 2022).'
Add-Run $p $t $false
$anchor = $p

# --- paragraph 5: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = [char]0x00C4 + 'ven Skogsstyrelsens egen v' + [char]0x00E4 + 'gledning f' + [char]0x00F6 + 'r h' + [char]0x00E4 + 'nsyn till kn' + [char]0x00E4 + 'rot ligger i linje med ovanst' + [char]0x00E5 + 'ende forskningsstudier. Av v' + [char]0x00E4 + 'gledningen framg' + [char]0x00E5 + 'r det att f' + [char]0x00F6 + 'r med h' + [char]0x00F6 + 'g sannolikhet kunna bevara befintliga f' + [char]0x00F6 + 'rekomster kr' + [char]0x00E4 + 'vs relativt stora avs' + [char]0x00E4 + 'ttningar av uppvuxen skog med slutet och relativt t' + [char]0x00E4 + 'tt kronskikt. Som riktlinje kan kr' + [char]0x00E4 + 'vas ett avst' + [char]0x00E5 + 'nd p' + [char]0x00E5 + ' 50 meter in fr' + [char]0x00E5 + 'n brynet f' + [char]0x00F6 + 'r att vidmakth' + [char]0x00E5 + 'lla ett fungerande mikroklimat. Detta inneb' + [char]0x00E4 + 'r att frist' + [char]0x00E5 + 'ende h' + [char]0x00E4 + 'nsynsytor f' + [char]0x00F6 + 'r m' + [char]0x00E5 + 'nga arter (k' + [char]0x00E4 + 'rlv' + [char]0x00E4 + 'xter, lavar och mossor) kan beh' + [char]0x00F6 + 'va ha en area ' + [char]0x00F6 + 'verstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) f' + [char]0x00F6 + 'r att bibeh' + [char]0x00E5 + 'lla lokalklimatet. ' + [char]0x00C4 + 'ven ganska sm' + [char]0x00E5 + ' f' + [char]0x00F6 + 'r' + [char]0x00E4 + 'ndringar i form av f' + [char]0x00F6 + 'r' + [char]0x00E4 + 'ndrade ljus- och fuktighetsf' + [char]0x00F6 + 'rh' + [char]0x00E5 + 'llanden, till exempel till f' + [char]0x00F6 + 'ljd av gallring, kan leda till att arten f' + [char]0x00F6 + 'rsvinner till f' + [char]0x00F6 + 'ljd av konkurrens med mera ljuskr' + [char]0x00E4 + 'vande och snabbv' + [char]0x00E4 + 'xande arter (Skogsstyrelsen, 2022).'
Add-Run $p $t $false
$anchor = $p

# --- paragraph 6: style='Heading2' ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -3
$t = 'Referenser - kn' + [char]0x00E4 + 'rot'
Add-Run $p $t $false
$anchor = $p

# --- paragraph 7: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = 'de Graaf M & Roberts M.R., 2009. '
Add-Run $p $t $false
$t = 'Short-term response of the herbaceous layer within leave patches after harvest. '
Add-Run $p $t $true
$t = 'Forest Ecology and Management 257, 1014-1025'
Add-Run $p $t $false
$anchor = $p

# --- paragraph 8: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. '
Add-Run $p $t $false
$t = 'Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. '
Add-Run $p $t $true
$t = 'Ecological Applications, 22, 2049-2064 '
Add-Run $p $t $false
$anchor = $p

# --- paragraph 9: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = 'Koelmeijer, I. A., Ehrl' + [char]0x00E9 + 'n, J., J' + [char]0x00F6 + 'nsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. '
Add-Run $p $t $false
$t = 'Interactive effects of drought and edge exposure on old-growth forest understory species. '
Add-Run $p $t $true
$t = 'Landscape Ecology, 37, sid 1839-1853'
Add-Run $p $t $false
$anchor = $p

# --- paragraph 10: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = 'Rudolphi, J., J' + [char]0x00F6 + 'nsson, M. T., & Gustafsson, L., 2014. '
Add-Run $p $t $false
$t = 'Biological legacies buffer local species extinction after logging. '
Add-Run $p $t $true
$t = 'Journal of Applied Ecology. 51, 53-62.'
Add-Run $p $t $false
$anchor = $p

# --- paragraph 11: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = 'Skogsstyrelsen, 2022. '
Add-Run $p $t $false
$t = 'V' + [char]0x00E4 + 'gledning f' + [char]0x00F6 + 'r h' + [char]0x00E4 + 'nsyn till kn' + [char]0x00E4 + 'rot. '
Add-Run $p $t $true
$t = 'https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/'
Add-Run $p $t $false
$anchor = $p

# --- paragraph 12: style=None ---
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = -1
$t = 'SLU Artdatabanken, 2021. '
Add-Run $p $t $false
$t = 'Artfaktablad. Naturv' + [char]0x00E5 + 'rd ' + [char]0x2013 + ' artfakta. '
Add-Run $p $t $true
$t = 'SLU Artdatabanken, Uppsala '
Add-Run $p $t $false
$anchor = $p

# --- date change in first-page header ---
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(2)
$hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
